$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row corresponding to participant "Julieta Hernandez" (row 3),
# shifting all rows below it up by one.
$ws.Rows.Item(3).Delete()
